$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.0125
$ws.Range("E2").Value = 4.35
$ws.Range("F2").Value = 0.0397
$ws.Range("C3").Value = 0.0012
$ws.Range("E3").Value = 0.418
$ws.Range("F3").Value = 0.519
$ws.Range("C4").Value = 0.0682
$ws.Range("E4").Value = 23.7
$ws.Range("F4").Value = 0.00000438
$ws.Range("C5").Value = 0.00225
$ws.Range("E5").Value = 0.783
$ws.Range("F5").Value = 0.378
$ws.Range("C6").Value = 0.16
$ws.Range("E6").Value = 55.7
$ws.Range("F6").Value = 0.0000000000369
$ws.Range("C7").Value = 0.0413
$ws.Range("E7").Value = 14.3
$ws.Range("F7").Value = 0.000265
$ws.Range("C8").Value = 0.00389
$ws.Range("E8").Value = 1.35
$ws.Range("F8").Value = 0.248
$ws.Range("C9").Value = 1.24
$ws.Range("E9").Value = 5.2
$ws.Range("F9").Value = 0.0248
$ws.Range("C10").Value = 9.83
$ws.Range("E10").Value = 41.1
$ws.Range("F10").Value = 0.00000000523
$ws.Range("C11").Value = 12.2
$ws.Range("E11").Value = 51.1
$ws.Range("F11").Value = 0.000000000162
$ws.Range("C12").Value = 9.53
$ws.Range("E12").Value = 39.8
$ws.Range("F12").Value = 0.00000000831
$ws.Range("C13").Value = 11.7
$ws.Range("E13").Value = 48.7
$ws.Range("F13").Value = 0.000000000364
$ws.Range("C14").Value = 0.135
$ws.Range("E14").Value = 0.566
$ws.Range("F14").Value = 0.454
$ws.Range("C15").Value = 0.635
$ws.Range("E15").Value = 2.65
$ws.Range("F15").Value = 0.106
$ws.Range("C16").Value = 0.00000831
$ws.Range("E16").Value = 0.691
$ws.Range("F16").Value = 0.408
$ws.Range("C17").Value = 0.000275
$ws.Range("E17").Value = 22.8
$ws.Range("F17").Value = 0.00000625
$ws.Range("C18").Value = 0.000734
$ws.Range("E18").Value = 61.1
$ws.Range("F18").Value = 0.00000000000665
$ws.Range("C19").Value = 0.000277
$ws.Range("E19").Value = 23.1
$ws.Range("F19").Value = 0.00000573
$ws.Range("C20").Value = 0.00104
$ws.Range("E20").Value = 86.6
$ws.Range("F20").Value = 0.00000000000000428
$ws.Range("C21").Value = 0.0000372
$ws.Range("E21").Value = 3.09
$ws.Range("F21").Value = 0.0819
$ws.Range("C22").Value = 0.000103
$ws.Range("E22").Value = 8.55
$ws.Range("F22").Value = 0.00431
$ws.Range("C23").Value = 1.1
$ws.Range("E23").Value = 1.75
$ws.Range("F23").Value = 0.189
$ws.Range("C24").Value = 6.93
$ws.Range("E24").Value = 11.0
$ws.Range("F24").Value = 0.00127
$ws.Range("C25").Value = 0.00495
$ws.Range("E25").Value = 0.00787
$ws.Range("F25").Value = 0.93
$ws.Range("C26").Value = 0.000244
$ws.Range("E26").Value = 0.000388
$ws.Range("F26").Value = 0.984
$ws.Range("C27").Value = 7.69
$ws.Range("E27").Value = 12.2
$ws.Range("F27").Value = 0.000712
$ws.Range("C28").Value = 0.000529
$ws.Range("E28").Value = 0.000841
$ws.Range("F28").Value = 0.977
$ws.Range("C29").Value = 1.53
$ws.Range("E29").Value = 2.42
$ws.Range("F29").Value = 0.123
$ws.Range("C30").Value = 69.9
$ws.Range("E30").Value = 0.376
$ws.Range("F30").Value = 0.541
$ws.Range("C31").Value = 23300.0
$ws.Range("E31").Value = 125.0
$ws.Range("F31").Value = 0.000000000000000000362
$ws.Range("C32").Value = 927.0
$ws.Range("E32").Value = 4.99
$ws.Range("F32").Value = 0.0277
$ws.Range("C33").Value = 29.1
$ws.Range("E33").Value = 0.156
$ws.Range("F33").Value = 0.693
$ws.Range("C34").Value = 1930.0
$ws.Range("E34").Value = 10.4
$ws.Range("F34").Value = 0.00172
$ws.Range("C35").Value = 30.5
$ws.Range("E35").Value = 0.164
$ws.Range("F35").Value = 0.686
$ws.Range("C36").Value = 1370.0
$ws.Range("E36").Value = 7.39
$ws.Range("F36").Value = 0.00777
$ws.Range("C37").Value = 12.8
$ws.Range("E37").Value = 1.01
$ws.Range("F37").Value = 0.318
$ws.Range("C38").Value = 118.0
$ws.Range("E38").Value = 9.31
$ws.Range("F38").Value = 0.00295
$ws.Range("C39").Value = 0.0582
$ws.Range("E39").Value = 0.0046
$ws.Range("F39").Value = 0.946
$ws.Range("C40").Value = 25.2
$ws.Range("E40").Value = 1.99
$ws.Range("F40").Value = 0.162
$ws.Range("C41").Value = 2.47
$ws.Range("E41").Value = 0.195
$ws.Range("F41").Value = 0.66
$ws.Range("C42").Value = 116.0
$ws.Range("E42").Value = 9.15
$ws.Range("F42").Value = 0.00318
$ws.Range("C43").Value = 6.91
$ws.Range("E43").Value = 0.546
$ws.Range("F43").Value = 0.462
